$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("C2").Value = 69
$ws.Range("D2").Value = 89
$ws.Range("E2").Value = 69
$ws.Range("H2").Value = 103
$ws.Range("I2").Value = 113
$ws.Range("E3").Value = 141
$ws.Range("F3").Value = 132
$ws.Range("H3").Value = 151
$ws.Range("J3").Value = 221
$ws.Range("K3").Value = 212
$ws.Range("B6").Value = 365
$ws.Range("C6").Value = 461
$ws.Range("D6").Value = 397
$ws.Range("E6").Value = 453
$ws.Range("F6").Value = 501
$ws.Range("G6").Value = 428
$ws.Range("H6").Value = 432
$ws.Range("I6").Value = 487
$ws.Range("J6").Value = 403
$ws.Range("K6").Value = 491
$ws.Range("B7").Value = 489
$ws.Range("C7").Value = 613
$ws.Range("D7").Value = 624
$ws.Range("E7").Value = 675
$ws.Range("F7").Value = 728
$ws.Range("G7").Value = 655
$ws.Range("H7").Value = 701
$ws.Range("I7").Value = 815
$ws.Range("J7").Value = 761
$ws.Range("K7").Value = 863

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I6").Value = 7
$ws.Range("B7").Value = 19
$ws.Range("F8").Value = 47
$ws.Range("K27").Value = 17
$ws.Range("E28").Value = 31
$ws.Range("F28").Value = 56
$ws.Range("K28").Value = 61
$ws.Range("E29").Value = 8
$ws.Range("D30").Value = 11
$ws.Range("B32").Value = 15
$ws.Range("E32").Value = 62
$ws.Range("J32").Value = 43
$ws.Range("C35").Value = 12
$ws.Range("E35").Value = 6
$ws.Range("J35").Value = 6
$ws.Range("D36").Value = 35
$ws.Range("E36").Value = 37
$ws.Range("G36").Value = 26
$ws.Range("J36").Value = 43
$ws.Range("I43").Value = 6
$ws.Range("K47").Value = 21
$ws.Range("G50").Value = 15
$ws.Range("H53").Value = 96
$ws.Range("B63").Value = 7
$ws.Range("C65").Value = 21
$ws.Range("D65").Value = 23
$ws.Range("H69").Value = 3
$ws.Range("J70").Value = 14
$ws.Range("H74").Value = 15
$ws.Range("J77").Value = 34
$ws.Range("E80").Value = 11
$ws.Range("J85").Value = 13
$ws.Range("B98").Value = 489
$ws.Range("C98").Value = 613
$ws.Range("D98").Value = 624
$ws.Range("E98").Value = 675
$ws.Range("F98").Value = 728
$ws.Range("G98").Value = 655
$ws.Range("H98").Value = 701
$ws.Range("I98").Value = 815
$ws.Range("J98").Value = 761
$ws.Range("K98").Value = 863

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("B5").Value = 17
$ws.Range("B6").Value = 19

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("D2").Value = 1
$ws.Range("D6").Value = 11

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("F6").Value = 33
$ws.Range("F7").Value = 47

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("B6").Value = 13
$ws.Range("E6").Value = 49
$ws.Range("J6").Value = 28
$ws.Range("B7").Value = 15
$ws.Range("E7").Value = 62
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 9
$ws.Range("D6").Value = 20
$ws.Range("G6").Value = 14
$ws.Range("J6").Value = 24
$ws.Range("D7").Value = 35
$ws.Range("E7").Value = 37
$ws.Range("G7").Value = 26
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("G5").Value = 14
$ws.Range("G6").Value = 15

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("E4").Value = 5
$ws.Range("E5").Value = 11

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("F3").Value = 11
$ws.Range("K3").Value = 18
$ws.Range("E6").Value = 20
$ws.Range("E7").Value = 31
$ws.Range("F7").Value = 56
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("H2").Value = 13
$ws.Range("H7").Value = 96

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 1
$ws.Range("C5").Value = 10
$ws.Range("E5").Value = 5
$ws.Range("C6").Value = 12
$ws.Range("E6").Value = 6
$ws.Range("J6").Value = 6

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("C2").Value = 3
$ws.Range("D5").Value = 22
$ws.Range("C6").Value = 21
$ws.Range("D6").Value = 23

$ws = $wb.Worksheets.Item("River North")
$ws.Range("H5").Value = 11
$ws.Range("H6").Value = 15

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 8

$ws = $wb.Worksheets.Item("New City")
$ws.Range("B4").Value = 5
$ws.Range("B5").Value = 7

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K4").Value = 13
$ws.Range("K5").Value = 17

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J4").Value = 6
$ws.Range("J5").Value = 13

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I2").Value = 2
$ws.Range("I6").Value = 6

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K5").Value = 9
$ws.Range("K6").Value = 21

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J5").Value = 6
$ws.Range("J6").Value = 14

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I4").Value = 3
$ws.Range("I5").Value = 7

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("E3").Value = 2
$ws.Range("E5").Value = 3
